# Update cryptos list: refresh Price (D) and Volume(1h) (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.333.94'
$ws.Range("E2").Value = '  -0.71%  '

$ws.Range("D3").Value = '1.710.51'
$ws.Range("E3").Value = '  -0.81%  '

$ws.Range("E4").Value = '  +0.09%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '224.54'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.58%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5288'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -1.04%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.004'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +0.09%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.06683'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +1.39%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2657'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -0.65%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07697'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -0.13%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.503'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -2.36%  '

$ws.Range("D13").Value = '1.946.45'
$ws.Range("E13").Value = '  -0.78%  '

$ws.Range("D14").Value = '1.711.17'
$ws.Range("E14").Value = '  -0.72%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5851'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +0.37%  '

$ws.Range("D16").Value = '0.0₅8209'
$ws.Range("E16").Value = '  -0.85%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '67.93'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -0.09%  '

$ws.Range("D18").Value = '27.360.87'

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '222.65'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +1.47%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.004'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +0.06%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.638'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -2.00%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.46'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -1.51%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.009'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -1.26%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.004'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.09%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '144.43'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -0.92%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.688'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -2.57%  '

$ws.Range("E27").Value = '  -2.23%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.235'
$ws.Range("D28").ClearFormats()

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '16.24'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -2.02%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05337'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -3.67%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.291'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -1.09%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.463'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -2.91%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.427'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -0.64%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.635'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -1.68%  '

$ws.Range("E35").Value = '  +0.54%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9519'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -1.49%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.393'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -1.26%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5848'
$ws.Range("D38").ClearFormats()

$ws.Range("D39").Value = '1.149.34'
$ws.Range("E39").Value = '  +8.78%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01636'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -0.92%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.784'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -2.13%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.004'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +0.12%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8393'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -1.94%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '100.96'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -0.35%  '

$ws.Range("D45").Value = '1.853.46'
$ws.Range("E45").Value = '  -0.79%  '

$ws.Range("D46").Value = '0.0₈112'
$ws.Range("E46").Value = '  -3.29%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '57.49'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -2.58%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4565'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +2.43%  '

$ws.Range("E49").Value = '  -0.07%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.069'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -1.66%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05218'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -0.40%  '
